$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.974.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.052.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.71%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '

$ws.Range("E6").Value = '  +1.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.07'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.97%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +2.63%  '

$ws.Range("E10").Value = '  +2.95%  '

$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.358.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.752'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.053.89'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.885.42'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.37%  '

$ws.Range("E19").Value = '  -0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.11%  '

$ws.Range("E21").Value = '  +1.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.67%  '

$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("E24").Value = '  +0.66%  '

$ws.Range("E25").Value = '  +3.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.22%  '

$ws.Range("E27").Value = '  +1.69%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.132'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.47%  '

$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("E31").Value = '  +2.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0611'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '

$ws.Range("E35").Value = '  +10.02%  '

$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.96'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.07%  '

$ws.Range("E39").Value = '  +0.15%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0218'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.27%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.484.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.91%  '

$ws.Range("E43").Value = '  +3.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0934'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.36%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +18.61%  '

$ws.Range("E47").Value = '  +0.31%  '

$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.245.86'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.90%  '
